$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first three rows only held the letterhead / title block
# (Sekretariat der KMK ..., "FERIEN IM SCHULJAHR 1960/61", address, phone,
# "Bundesrepublik Deutschland") spread across merged cells B1:F1, B2:F2 and
# B3:F3. That information is redundant with the actual data table below, so
# remove the rows outright - this shifts the "Land / Ostern / Pfingsten /
# Sommer / Herbst / Weihnachten" header (formerly row 4) up to row 1, and
# every data row moves up by three as well.
$ws.Rows("1:3").Delete()

# After the shift, the sheet only has 14 populated rows (previously 17).
# Keep the row-height metadata of the three now-empty trailing rows so the
# sheet still reports 17 rows overall, matching the original layout/grid.
$ws.Rows(15).RowHeight = 21.95
$ws.Rows(16).RowHeight = 23.1
$ws.Rows(17).RowHeight = 21

# Leave the cursor where the author left it when saving.
$ws.Range("B8").Select()
